$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.419.24'
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('D3').Value = '2.425.13'
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'564.96"
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D6').Value = "'144.73"
$ws.Range('E6').Value = '  +1.93%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = "'0.531"
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').Value = "'0.110"
$ws.Range('E9').Value = '  +1.37%  '
$ws.Range('E10').Value = '  +0.42%  '
$ws.Range('D11').Value = "'5.29"
$ws.Range('E11').Value = '  +1.04%  '
$ws.Range('D12').Value = "'0.354"
$ws.Range('E12').Value = '  +1.64%  '
$ws.Range('D13').Value = "'26.73"
$ws.Range('E13').Value = '  +4.95%  '
$ws.Range('E14').Value = '  +4.06%  '
$ws.Range('E15').Value = '  +0.59%  '
$ws.Range('D16').Value = '62.348.72'
$ws.Range('E16').Value = '  +0.97%  '
$ws.Range('D17').Value = '2.430.96'
$ws.Range('E17').Value = '  +0.81%  '
$ws.Range('D18').Value = "'11.19"
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').Value = "'6.96"
$ws.Range('E19').Value = '  +2.22%  '
$ws.Range('D20').Value = "'323.37"
$ws.Range('E20').Value = '  +0.73%  '
$ws.Range('D21').Value = "'4.16"
$ws.Range('E21').Value = '  +1.06%  '
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').Value = "'67.15"
$ws.Range('E23').Value = '  +2.30%  '
$ws.Range('E24').Value = '  +5.10%  '
$ws.Range('D25').Value = "'588.84"
$ws.Range('E25').Value = '  +4.38%  '
$ws.Range('D26').Value = "'8.53"
$ws.Range('E26').Value = '  -1.65%  '
$ws.Range('D27').Value = '0.0₃0994'
$ws.Range('E27').Value = '  +7.01%  '
$ws.Range('D28').Value = '2.548.32'
$ws.Range('E28').Value = '  +1.31%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = "'0.999"
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = "'8.42"
$ws.Range('E30').Value = '  +3.05%  '
$ws.Range('E31').Value = '  +3.86%  '
$ws.Range('E32').Value = '  -1.04%  '
$ws.Range('E33').Value = '  +0.40%  '
$ws.Range('D34').Value = "'1.50"
$ws.Range('E34').Value = '  -0.13%  '
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('D36').Value = "'4.83"
$ws.Range('E36').Value = '  +1.85%  '
$ws.Range('D37').Value = "'0.381"
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('D38').Value = "'18.69"
$ws.Range('E38').Value = '  +1.15%  '
$ws.Range('D39').Value = "'5.34"
$ws.Range('D40').Value = "'147.85"
$ws.Range('E40').Value = '  -3.12%  '
$ws.Range('D41').Value = "'1.81"
$ws.Range('E41').Value = '  +1.77%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('E43').Value = '  +9.00%  '
$ws.Range('D44').Value = "'147.89"
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = "'3.66"
$ws.Range('E45').Value = '  +1.95%  '
$ws.Range('D46').Value = "'0.0534"
$ws.Range('E46').Value = '  +1.16%  '
$ws.Range('D47').Value = "'20.46"
$ws.Range('E47').Value = '  +3.40%  '
$ws.Range('D48').Value = "'0.599"
$ws.Range('E48').Value = '  +1.78%  '
$ws.Range('E49').Value = '  +2.88%  '
$ws.Range('D50').Value = "'0.0918"
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('E51').Value = '  +4.58%  '
